$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update individual cell values in rows 2-25 (row numbers unaffected by later row deletions) ---

# Row 2
$ws.Range("D2").Value = -13.5
$ws.Range("F2").Value = 18.03

# Row 3
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()

# Row 4
$ws.Range("E4").Value = -6.4

# Row 5
$ws.Range("E5").ClearContents()

# Row 6
$ws.Range("D6").ClearContents()

# Row 8
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()

# Row 12
$ws.Range("D12").Value = -14.1

# Row 13
$ws.Range("F13").Value = 17.1

# Row 14
$ws.Range("D14").ClearContents()

# Row 15
$ws.Range("E15").Value = -8.4

# Row 18
$ws.Range("E18").Value = -8.5

# Row 19
$ws.Range("E19").ClearContents()
$ws.Range("F19").ClearContents()

# Row 20
$ws.Range("D20").Value = -14

# Row 21
$ws.Range("D21").Value = -14.3

# Row 22
$ws.Range("E22").ClearContents()

# Row 23
$ws.Range("D23").ClearContents()
$ws.Range("E23").Value = -7

# Row 24
$ws.Range("D24").ClearContents()

# Row 25
$ws.Range("E25").Value = -7.1
$ws.Range("F25").Value = 16.6

# --- Remove the "RM 232" row (row 26) and the "SC 92" row (originally row 28, now row 27) ---
$ws.Rows("26").Delete()
$ws.Rows("27").Delete()

# --- Update values in the rows that shifted up after the deletions ---

# Row 26 (was "SC 5")
$ws.Range("C26").Value = 10.8

# Row 27 (was "SC 101")
$ws.Range("C27").ClearContents()
$ws.Range("E27").ClearContents()

# Row 28 (was "SC 105")
$ws.Range("F28").Value = 17.44

# Row 30 (was "SC 120")
$ws.Range("C30").Value = 11.4

# Row 31 (was "SC 132")
$ws.Range("D31").Value = -13.7
$ws.Range("F31").ClearContents()

# Row 32 (was "SC 193")
$ws.Range("C32").ClearContents()
$ws.Range("F32").Value = 17.39

# Row 33 (was "SC 232")
$ws.Range("D33").Value = -14.1
